$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "REQUISITION:P2318L01H0"
$ws.Range("C2").Value = "Erez A Minka"
$ws.Range("E2").Value = "Wichita Falls,"

# Zip code looks numeric - force text entry so it stays a shared string,
# then restore the cell's style so no extra number format sticks around.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "76308"
$ws.Range("F2").Style = "Normal"

$ws.Range("G2").Value = "940-687-3376"
$ws.Range("H2").Value = "4327 Barnett Road Wichita Falls, TX 763102303"
$ws.Range("I2").Value = "Askins, Sammie"
$ws.Range("J2").Value = "(M/79)"

# Date of birth looks like a date - same text-forcing trick.
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "1944-06-13"
$ws.Range("K2").Style = "Normal"

$ws.Range("L2").Value = "100 Chaparral Drive"
$ws.Range("M2").Value = "SAAS0001"
$ws.Range("N2").Value = "SAAS0001"
$ws.Range("O2").Value = "Medicare = Texas"
$ws.Range("P2").Value = "(7WJ3UD1AE99)"
$ws.Range("Q2").Value = "Blue Shield = Texas Askins,"
$ws.Range("T2").Value = "Anterior Left Upper Arm - Central"
$ws.Range("U2").Value = "Neoplasm of uncertain behavior of skin"
$ws.Range("V2").Value = "Biopsy (Tangential (Shave))"
$ws.Range("W2").Value = "Pink papule (Anterior Left Upper Arm = Central)"
$ws.Range("X2").Value = "Submandibular Neck - Left"
$ws.Range("Y2").Value = "Neoplasm of uncertain behavior of skin"
$ws.Range("Z2").Value = "Biopsy (Tangential (Shave))"
$ws.Range("AA2").Value = "Exam: Pink papule (Submandibular Neck = Left)"
$ws.Range("AB2").Value = "KY"

# RequisitionNumber zip style numeric text again.
$ws.Range("AC2").NumberFormat = "@"
$ws.Range("AC2").Value = "405124601"
$ws.Range("AC2").Style = "Normal"
